$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PFAS measurement values (rows 22-31 and 33), matching the
# revised contamination profile figures.
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = 7.9

$ws.Range("B23").Value = 8.5
$ws.Range("C23").Value = 8.1

$ws.Range("B24").Value = 3

$ws.Range("B25").Value = 8.199999999999999
$ws.Range("C25").Value = 5.8

$ws.Range("B26").Value = 1.6
$ws.Range("C26").Value = 2

$ws.Range("B27").Value = 6.1

$ws.Range("B28").Value = 12.1
$ws.Range("C28").Value = 3.6

$ws.Range("B29").Value = 2.4
$ws.Range("C29").Value = 1.6

$ws.Range("B30").Value = 44.5
$ws.Range("C30").Value = 58.6

$ws.Range("B31").Value = 6.3
$ws.Range("C31").Value = 9.199999999999999

$ws.Range("B33").Value = 1.2

# Remove the HBCDD rows (a-HBCDD, b-HBCDD, g-HBCDD) entirely, shifting
# everything below them upward.
$ws.Range("A34:C36").EntireRow.Delete()
